# "fixed export and fixing maps"
# - rename the sheet from the generic "1" to the municipality name
# - drop the old 1989/2002 census-comparison columns, keeping only the
#   2014 figure (which slides from column D into column B)
# - clear out the leftover "(census results)" note that used to sit
#   under the title
# - drop the now-empty spacer row that used to separate the note from
#   the "(sq. km)" unit label

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "წალენჯიხა"

# Remove the 1989 and 2002 data columns (old B:C); the 2014 column
# (old D) shifts left into column B.
$ws.Range("B:C").Delete()

# The census-results note under the title is no longer used.
$ws.Range("A2").ClearContents()

# The blank spacer row (old row 3) collapses away entirely.
$ws.Range("3:3").Delete()

$ws.Range("A2").Select()
